$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1651
$ws.Range("I6").Value = 2722
$ws.Range("J6").Value = 580
$ws.Range("K6").Value = 8166
$ws.Range("L6").Value = 1740
$ws.Range("M6").Value = -8054
$ws.Range("N6").Value = -1964
$ws.Range("H9").Value = 266.4375
$ws.Range("I9").Value = 321.4
$ws.Range("K9").Value = 321.4
$ws.Range("M9").Value = -152.4
$ws.Range("H31").Value = 3483.3333
$ws.Range("I31").Value = 1850
$ws.Range("K31").Value = 5550
$ws.Range("M31").Value = -5320
$ws.Range("H74").Value = 6658.8276
$ws.Range("I74").Value = 7076.28
$ws.Range("J74").Value = 4049.75
$ws.Range("K74").Value = 7076.28
$ws.Range("L74").Value = 4049.75
$ws.Range("M74").Value = -6140.28
$ws.Range("N74").Value = -5921.75
$ws.Range("H77").Value = 6658.8276
$ws.Range("I77").Value = 7076.28
$ws.Range("J77").Value = 4049.75
$ws.Range("K77").Value = 35381.4
$ws.Range("L77").Value = 20248.75
$ws.Range("M77").Value = -30701.4
$ws.Range("N77").Value = -29608.75
$ws.Range("H88").Value = 1766.6154
$ws.Range("I88").Value = 2803
$ws.Range("J88").Value = 1118.875
$ws.Range("K88").Value = 2803
$ws.Range("L88").Value = 1118.875
$ws.Range("M88").Value = -2397
$ws.Range("N88").Value = -1930.875
$ws.Range("H91").Value = 1766.6154
$ws.Range("I91").Value = 2803
$ws.Range("J91").Value = 1118.875
$ws.Range("K91").Value = 2803
$ws.Range("L91").Value = 1118.875
$ws.Range("M91").Value = -1399
$ws.Range("N91").Value = -3926.875
$ws.Range("H100").Value = 3153.5
$ws.Range("I100").Value = 3147.8572
$ws.Range("K100").Value = 3147.8572
$ws.Range("M100").Value = -2606.8572
$ws.Range("H106").Value = 2107.8235
$ws.Range("I106").Value = 1055.6666
$ws.Range("K106").Value = 1055.6666
$ws.Range("M106").Value = -424.6666
$ws.Range("H129").Value = 626341.9
$ws.Range("I129").Value = 500957.1
$ws.Range("J129").Value = 835316.5600000001
$ws.Range("K129").Value = 1502871.3
$ws.Range("L129").Value = 2505949.68
$ws.Range("M129").Value = -1497871.3
$ws.Range("N129").Value = -2515949.68

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 487867.7
$ws.Range("J32").Value = 2186.75
$ws.Range("L32").Value = 2186.75
$ws.Range("N32").Value = -2760.75
$ws.Range("H45").Value = 2023.8889
$ws.Range("I45").Value = 2023.8889
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2023.8889
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1646.8889
$ws.Range("H132").Value = 3494.889
$ws.Range("I132").Value = 1922.0476
$ws.Range("K132").Value = 5766.142800000001
$ws.Range("M132").Value = -3236.142800000001
$ws.Range("H140").Value = 100129.29
$ws.Range("J140").Value = 50150.832
$ws.Range("L140").Value = 50150.832
$ws.Range("N140").Value = -60510.832

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2811.5
$ws.Range("I86").Value = 2371
$ws.Range("J86").Value = 4133
$ws.Range("K86").Value = 2371
$ws.Range("L86").Value = 4133
$ws.Range("M86").Value = -1248
$ws.Range("N86").Value = -6379
$ws.Range("H89").Value = 2811.5
$ws.Range("I89").Value = 2371
$ws.Range("J89").Value = 4133
$ws.Range("K89").Value = 11855
$ws.Range("L89").Value = 20665
$ws.Range("M89").Value = -6239
$ws.Range("N89").Value = -31897
$ws.Range("H105").Value = 2025.1
$ws.Range("I105").Value = 1608.1428
$ws.Range("K105").Value = 1608.1428
$ws.Range("M105").Value = 138.8571999999999
$ws.Range("H107").Value = 1617.8334
$ws.Range("I107").Value = 1741.4
$ws.Range("K107").Value = 1741.4
$ws.Range("M107").Value = 178.5999999999999
$ws.Range("H134").Value = 4537681
$ws.Range("I134").Value = 3790338.8
$ws.Range("K134").Value = 11371016.4
$ws.Range("M134").Value = -11368481.4

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3767.2
$ws.Range("I105").Value = 4130.5
$ws.Range("K105").Value = 4130.5
$ws.Range("M105").Value = -2383.5
$ws.Range("H134").Value = 5101.657
$ws.Range("I134").Value = 3642.389
$ws.Range("K134").Value = 10927.167
$ws.Range("M134").Value = -8392.167000000001

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7190.9287
$ws.Range("J3").Value = 15000
$ws.Range("L3").Value = 45000
$ws.Range("N3").Value = -45224
$ws.Range("H7").Value = 476.14285
$ws.Range("I7").Value = 966.6667
$ws.Range("J7").Value = 108.25
$ws.Range("K7").Value = 2900.0001
$ws.Range("L7").Value = 324.75
$ws.Range("M7").Value = -2788.0001
$ws.Range("N7").Value = -548.75
$ws.Range("H17").Value = 726.7778
$ws.Range("I17").Value = 1014.8333
$ws.Range("J17").Value = 150.66667
$ws.Range("K17").Value = 3044.4999
$ws.Range("L17").Value = 452.00001
$ws.Range("M17").Value = -2875.4999
$ws.Range("N17").Value = -790.00001
$ws.Range("H131").Value = 3801.35
$ws.Range("J131").Value = 6081.727
$ws.Range("L131").Value = 18245.181
$ws.Range("N131").Value = -28325.181

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1076.5834
$ws.Range("I97").Value = 960.3
$ws.Range("J97").Value = 1159.6428
$ws.Range("K97").Value = 960.3
$ws.Range("L97").Value = 1159.6428
$ws.Range("M97").Value = -464.3
$ws.Range("N97").Value = -2151.6428
$ws.Range("H102").Value = 2366.4167
$ws.Range("I102").Value = 2408.818
$ws.Range("K102").Value = 2408.818
$ws.Range("M102").Value = -786.8180000000002

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3827
$ws.Range("I7").Value = 2908.75
$ws.Range("K7").Value = 2908.75
$ws.Range("M7").Value = -2796.75
$ws.Range("H16").Value = 1658.6154
$ws.Range("I16").Value = 1584.5
$ws.Range("K16").Value = 1584.5
$ws.Range("M16").Value = -1414.5
$ws.Range("H21").Value = 7500
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 7500
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 7500
$ws.Range("H22").Value = 2374.7585
$ws.Range("J22").Value = 3297.3333
$ws.Range("L22").Value = 3297.3333
$ws.Range("N22").Value = -3887.3333
$ws.Range("H27").Value = 2374.7585
$ws.Range("J27").Value = 3297.3333
$ws.Range("L27").Value = 3297.3333
$ws.Range("N27").Value = -3511.3333
$ws.Range("H68").Value = 5796.5
$ws.Range("I68").Value = 3808.25
$ws.Range("J68").Value = 13749.5
$ws.Range("K68").Value = 3808.25
$ws.Range("L68").Value = 13749.5
$ws.Range("M68").Value = -3059.25
$ws.Range("N68").Value = -15247.5
$ws.Range("H71").Value = 5796.5
$ws.Range("I71").Value = 3808.25
$ws.Range("J71").Value = 13749.5
$ws.Range("K71").Value = 19041.25
$ws.Range("L71").Value = 68747.5
$ws.Range("M71").Value = -15297.25
$ws.Range("N71").Value = -76235.5
$ws.Range("H93").Value = 4632.1665
$ws.Range("I93").Value = 1250
$ws.Range("K93").Value = 1250
$ws.Range("M93").Value = -2
$ws.Range("H122").Value = 3766.3333
$ws.Range("I122").Value = 3619.6
$ws.Range("K122").Value = 10858.8
$ws.Range("M122").Value = -8408.799999999999
$ws.Range("H126").Value = 3827
$ws.Range("I126").Value = 2908.75
$ws.Range("K126").Value = 8726.25
$ws.Range("M126").Value = -6256.25

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 4896.826
$ws.Range("I13").Value = 4934.619
$ws.Range("J13").Value = 4500
$ws.Range("K13").Value = 4934.619
$ws.Range("L13").Value = 4500
$ws.Range("M13").Value = -4794.619
$ws.Range("N13").Value = -4780
$ws.Range("H122").Value = 77390.336
$ws.Range("J122").Value = 115006.5
$ws.Range("L122").Value = 345019.5
$ws.Range("N122").Value = -349919.5
$ws.Range("H126").Value = 1964.375
$ws.Range("I126").Value = 2108.182
$ws.Range("J126").Value = 1648
$ws.Range("K126").Value = 6324.545999999999
$ws.Range("L126").Value = 4944
$ws.Range("M126").Value = -3854.545999999999
$ws.Range("N126").Value = -9884
$ws.Range("H132").Value = 3088499.5
$ws.Range("I132").Value = 3335323.5
$ws.Range("K132").Value = 10005970.5
$ws.Range("M132").Value = -10003440.5

# ===== Special structural changes =====
$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("N45").ClearContents()

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("M21").ClearContents()
$wsLTW.Range("N21").Value = -7848
